{"js": "// Remove the trailing \" [Link]\" placeholder text that follows the\n// \"Distributed System & Database\" (Paxos) bullet at the very end of the\n// resume. Three other project bullets end with a similar-looking\n// \"[Link]\" marker, but those are real hyperlinks built from separate\n// \" [\" + hyperlink(\"Link\") + \"]\" runs; only this last paragraph has a\n// plain, unlinked \" [Link]\" text run, so scope the search to that one\n// paragraph instead of searching the whole body.\nconst body = context.document.body;\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\nconst results = lastParagraph.search(\" [Link]\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  for (const r of results.items) {\n    r.delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \" [Link]\" placeholder text that follows the\n# \"Distributed System & Database\" (Paxos) bullet at the very end of the\n# resume. Three other project bullets end with a similar-looking\n# \"[Link]\" marker, but those are real hyperlinks built from separate\n# \" [\" + hyperlink(\"Link\") + \"]\" runs; only the LAST paragraph in the\n# document has a plain, unlinked \" [Link]\" text run, so scope the Find\n# to that paragraph's Range instead of searching the whole document.\n$d = $word.ActiveDocument\n\n$lastParagraph = $d.Paragraphs.Last\n$rng = $lastParagraph.Range\n\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \" [Link]\"\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWholeWord = $false\n$rng.Find.MatchWildcards = $false\n$rng.Find.Forward = $true\n\n$found = $rng.Find.Execute()\n\nif ($found) {\n    $rng.Delete()\n}\n"}
